$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.003208871385164791, 0.0000005461030343489881, 3.537761648806719, 10.19245300693656, 13.73342407323148)
    3  = @(0.1190320826869504, 0.002571899574220771, 0.7527432677738641, 0.4942365360607697, 1.368583786095805)
    4  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    5  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    6  = @(0.01293466051926884, 0.306821227259698, 3.537761648806719, 10.19245300693656, 14.04997054352224)
    7  = @(0.01293466051926884, 0.002571899574220771, 0.7527432677738641, 0.4942365360607697, 1.262486363928123)
    8  = @(1.455362044514542, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 25.99591228164478)
    9  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    10 = @(0.04271373187048222, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.330410019770453)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
